$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2782779633998871
$ws.Range("B1").Value = 0.1977447122335434
$ws.Range("C1").Value = 0.2158285528421402
$ws.Range("D1").Value = 4.53001880645752
$ws.Range("E1").Value = 1.816592216491699
